# Sprint 2 Retrospective Presentation - slide 2 text tweak
# "Make sure you use it – Fridge expiry trackers and push notifications"
#   -> "Make sure you use it – Fridge expiry trackers and notifications"

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)          # "Content Placeholder 2"
$tr = $sh.TextFrame.TextRange

# Third bullet paragraph (level 1): "Make sure you use it – ..."
$para = $tr.Paragraphs(3, 1)
$run  = $para.Runs(1, 1)

$run.Text = "Make sure you use it " + [char]0x2013 + " Fridge expiry trackers and notifications"
